$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "5ª persona():" paragraph (the one with nothing between
# "persona(" and "):") and rewrite its text so the whole run (and the
# proofErr markers wrapping "persona(") are cleared out in one shot.
# ------------------------------------------------------------------
$paras = $d.Paragraphs
$targetPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t -like "5*persona(*):*") {
        $targetPara = $p
        break
    }
}

$newSentence = "5ª persona(Alberto):"

$pRange = $targetPara.Range
# Exclude the trailing paragraph mark from the replacement range.
$body = $d.Range($pRange.Start, $pRange.End - 1)
$paraStart = $body.Start
$body.Text = $newSentence

# ------------------------------------------------------------------
# Re-split "persona(" and "Alberto" into their own runs (mirrors what
# Word does when text is typed/edited in place) by toggling a format
# property on, then back off, for each sub-range.
# ------------------------------------------------------------------
$prefix = "5ª "
$personaWord = "persona("
$nameWord = "Alberto"

$personaStart = $paraStart + $prefix.Length
$personaEnd = $personaStart + $personaWord.Length
$nameStart = $personaEnd
$nameEnd = $nameStart + $nameWord.Length

$personaRng = $d.Range($personaStart, $personaEnd)
$personaRng.Bold = 1
$personaRng.Bold = 0

$nameRng = $d.Range($nameStart, $nameEnd)
$nameRng.Bold = 1
$nameRng.Bold = 0

# ------------------------------------------------------------------
# Move the (singleton) "_GoBack" bookmark so it now sits right after
# the newly-typed name, before the closing "):" — re-adding it here
# automatically removes it from its old location (end of the
# "4ª persona(José)" paragraph).
# ------------------------------------------------------------------
$bmRange = $d.Range($nameEnd, $nameEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
